# 自动更新Excel文件 - 2025-10-18 23:11:15
# For each data row, decrement the "remaining days" (column E) by 1.
# When the remaining days would drop to 0 (i.e. a row's E value was 1),
# reset E back to the total days (column D) and roll the start date
# (column F, stored as an integer YYYYMMDD) forward by D days.
# Rows whose start date cannot be parsed as a valid YYYYMMDD date are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($dVal -eq $null -or $eVal -eq $null -or $fVal -eq $null) {
        continue
    }
    if ($dVal -eq "" -or $eVal -eq "" -or $fVal -eq "") {
        continue
    }

    $totalDays = [int]$dVal
    $remaining = [int]$eVal
    $startStr = [string][int64]$fVal

    if ($startStr.Length -ne 8) {
        # malformed date, e.g. "202510929" - skip this row
        continue
    }

    $year = [int]$startStr.Substring(0, 4)
    $month = [int]$startStr.Substring(4, 2)
    $day = [int]$startStr.Substring(6, 2)

    $validDate = $true
    try {
        $startDate = Get-Date -Year $year -Month $month -Day $day
    } catch {
        $validDate = $false
    }

    if (-not $validDate) {
        continue
    }

    if ($remaining -eq 1) {
        $newStartDate = $startDate.AddDays($totalDays)
        $newF = [int]($newStartDate.ToString("yyyyMMdd"))
        $eCell.Value2 = $totalDays
        $fCell.Value2 = $newF
    } else {
        $eCell.Value2 = $remaining - 1
    }
}
